# "Version 2." -> "Version 1."
#
# Original run layout inside the single paragraph:
#   [0,5)  "Versi"
#   [5,7)  "on"
#   [7,9)  " 2"
#   (bookmarkStart/bookmarkEnd "_GoBack")
#   [9,10) "."
#
# Target run layout:
#   [0,7)  "Version"          (merge of "Versi" + "on")
#   [7,10) " 1."              (merge of " 2" -> " 1" and the trailing ".")
#   (bookmarkStart/bookmarkEnd "_GoBack", now trailing in the paragraph)

$d = $word.ActiveDocument

# --- Merge "Versi" + "on" into a single run reading "Version" ---
# Word only splits/merges runs when a real text change happens, so first
# write a temporary, differing value, then correct it back to "Version".
$rVersion = $d.Range(0, 7)
$rVersion.Text = "Versionx"
$rVersionFix = $d.Range(0, 8)
$rVersionFix.Text = "Version"

# --- Change the version number "2" -> "1" ---
$rNum = $d.Range(8, 9)
$rNum.Text = "1"

# --- Fold the trailing "." run into the " 1" run, producing " 1." ---
# (this keeps the bookmark, which sits between them, intact and trailing)
$rOne = $d.Range(7, 9)
$rOne.InsertAfter(".")

# --- Remove the now-duplicated old "." run ---
$rOldDot = $d.Range(10, 11)
$rOldDot.Delete()
